$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pdgfb"
$ws.Range("C2").Value = "Pdgfrb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 29.253501
$ws.Range("H2").Value = 87.760503
$ws.Range("I2").Value = 0.7876335333413836
$ws.Range("J2").Value = 0.7876335333413838
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 4.389289333333333
$ws.Range("N2").Value = 13.167868
$ws.Range("O2").Value = 0.03995844413671427
$ws.Range("P2").Value = 0.03995844413671427
$ws.Range("Q2").Value = 128.402079901956
$ws.Range("R2").Value = 1155.618719117604
$ws.Range("S2").Value = 0.03147261054222456
$ws.Range("T2").Value = 0.03147261054222456

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pdgfb"
$ws.Range("C3").Value = "Pdgfrb"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 29.253501
$ws.Range("H3").Value = 87.760503
$ws.Range("I3").Value = 0.7876335333413836
$ws.Range("J3").Value = 0.7876335333413838
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 57.89762366666667
$ws.Range("N3").Value = 173.692871
$ws.Range("O3").Value = 0.5270782546422108
$ws.Range("P3").Value = 0.5270782546422108
$ws.Range("Q3").Value = 1693.708191830457
$ws.Range("R3").Value = 15243.37372647411
$ws.Range("S3").Value = 0.415144508051254
$ws.Range("T3").Value = 0.4151445080512541

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pdgfb"
$ws.Range("C4").Value = "Pdgfrb"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 29.253501
$ws.Range("H4").Value = 87.760503
$ws.Range("I4").Value = 0.7876335333413836
$ws.Range("J4").Value = 0.7876335333413838
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.129282
$ws.Range("N4").Value = 0.387846
$ws.Range("O4").Value = 0.001176934848120294
$ws.Range("P4").Value = 0.001176934848120294
$ws.Range("Q4").Value = 3.781951116282
$ws.Range("R4").Value = 34.037560046538
$ws.Range("S4").Value = 0.0009269933529375922
$ws.Range("T4").Value = 0.0009269933529375923

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Pdgfb"
$ws.Range("C5").Value = "Pdgfrb"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 29.253501
$ws.Range("H5").Value = 87.760503
$ws.Range("I5").Value = 0.7876335333413836
$ws.Range("J5").Value = 0.7876335333413838
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 47.43015733333333
$ws.Range("N5").Value = 142.290472
$ws.Range("O5").Value = 0.4317863663729547
$ws.Range("P5").Value = 0.4317863663729548
$ws.Range("Q5").Value = 1387.498154980824
$ws.Range("R5").Value = 12487.48339482742
$ws.Range("S5").Value = 0.3400894213949675
$ws.Range("T5").Value = 0.3400894213949676

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Pdgfb"
$ws.Range("C6").Value = "Pdgfrb"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.749137666666667
$ws.Range("H6").Value = 14.247413
$ws.Range("I6").Value = 0.1278677748937237
$ws.Range("J6").Value = 0.1278677748937237
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 4.389289333333333
$ws.Range("N6").Value = 13.167868
$ws.Range("O6").Value = 0.03995844413671427
$ws.Range("P6").Value = 0.03995844413671427
$ws.Range("Q6").Value = 20.84533930283155
$ws.Range("R6").Value = 187.608053725484
$ws.Range("S6").Value = 0.005109397339976814
$ws.Range("T6").Value = 0.005109397339976816

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Pdgfb"
$ws.Range("C7").Value = "Pdgfrb"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.749137666666667
$ws.Range("H7").Value = 14.247413
$ws.Range("I7").Value = 0.1278677748937237
$ws.Range("J7").Value = 0.1278677748937237
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 57.89762366666667
$ws.Range("N7").Value = 173.692871
$ws.Range("O7").Value = 0.5270782546422108
$ws.Range("P7").Value = 0.5270782546422108
$ws.Range("Q7").Value = 274.9637853658581
$ws.Range("R7").Value = 2474.674068292723
$ws.Range("S7").Value = 0.067396323615967
$ws.Range("T7").Value = 0.06739632361596701

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Pdgfb"
$ws.Range("C8").Value = "Pdgfrb"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.749137666666667
$ws.Range("H8").Value = 14.247413
$ws.Range("I8").Value = 0.1278677748937237
$ws.Range("J8").Value = 0.1278677748937237
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.129282
$ws.Range("N8").Value = 0.387846
$ws.Range("O8").Value = 0.001176934848120294
$ws.Range("P8").Value = 0.001176934848120294
$ws.Range("Q8").Value = 0.6139780158220001
$ws.Range("R8").Value = 5.525802142398001
$ws.Range("S8").Value = 0.0001504920402240247
$ws.Range("T8").Value = 0.0001504920402240247

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Pdgfb"
$ws.Range("C9").Value = "Pdgfrb"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.749137666666667
$ws.Range("H9").Value = 14.247413
$ws.Range("I9").Value = 0.1278677748937237
$ws.Range("J9").Value = 0.1278677748937237
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 47.43015733333333
$ws.Range("N9").Value = 142.290472
$ws.Range("O9").Value = 0.4317863663729547
$ws.Range("P9").Value = 0.4317863663729548
$ws.Range("Q9").Value = 225.2523467276596
$ws.Range("R9").Value = 2027.271120548936
$ws.Range("S9").Value = 0.05521156189755588
$ws.Range("T9").Value = 0.0552115618975559

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Pdgfb"
$ws.Range("C10").Value = "Pdgfrb"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.138366333333334
$ws.Range("H10").Value = 9.415099000000001
$ws.Range("I10").Value = 0.08449869176489255
$ws.Range("J10").Value = 0.08449869176489258
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 4.389289333333333
$ws.Range("N10").Value = 13.167868
$ws.Range("O10").Value = 0.03995844413671427
$ws.Range("P10").Value = 0.03995844413671427
$ws.Range("Q10").Value = 13.77519787099245
$ws.Range("R10").Value = 123.976780838932
$ws.Range("S10").Value = 0.003376436254512898
$ws.Range("T10").Value = 0.003376436254512899

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Pdgfb"
$ws.Range("C11").Value = "Pdgfrb"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.138366333333334
$ws.Range("H11").Value = 9.415099000000001
$ws.Range("I11").Value = 0.08449869176489255
$ws.Range("J11").Value = 0.08449869176489258
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 57.89762366666667
$ws.Range("N11").Value = 173.692871
$ws.Range("O11").Value = 0.5270782546422108
$ws.Range("P11").Value = 0.5270782546422108
$ws.Range("Q11").Value = 181.7039528954699
$ws.Range("R11").Value = 1635.335576059229
$ws.Range("S11").Value = 0.04453742297498972
$ws.Range("T11").Value = 0.04453742297498973

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Pdgfb"
$ws.Range("C12").Value = "Pdgfrb"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.138366333333334
$ws.Range("H12").Value = 9.415099000000001
$ws.Range("I12").Value = 0.08449869176489255
$ws.Range("J12").Value = 0.08449869176489258
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.129282
$ws.Range("N12").Value = 0.387846
$ws.Range("O12").Value = 0.001176934848120294
$ws.Range("P12").Value = 0.001176934848120294
$ws.Range("Q12").Value = 0.4057342763060001
$ws.Range("R12").Value = 3.651608486754001
$ws.Range("S12").Value = 0.00009944945495867738
$ws.Range("T12").Value = 0.00009944945495867742

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Pdgfb"
$ws.Range("C13").Value = "Pdgfrb"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.138366333333334
$ws.Range("H13").Value = 9.415099000000001
$ws.Range("I13").Value = 0.08449869176489255
$ws.Range("J13").Value = 0.08449869176489258
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 47.43015733333333
$ws.Range("N13").Value = 142.290472
$ws.Range("O13").Value = 0.4317863663729547
$ws.Range("P13").Value = 0.4317863663729548
$ws.Range("Q13").Value = 148.8532089596365
$ws.Range("R13").Value = 1339.678880636728
$ws.Range("S13").Value = 0.03648538308043127
$ws.Range("T13").Value = 0.03648538308043128
